$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update invoice date
$ws.Range("H8").Value2 = "06 - 03 - 2020"

# Update the invoice description (scope of work)
$ws.Range("C14").Value2 = "Inventory Management System (IMS) as GUI Desktop Program with below capabilities:`n- Login Screen for Admin Access`n- Full customer/supplier invoice entry`n- Import customer/supplier invoices`n- Import customer/supplier/product list records`n- Re-order levels`n- Report on Stocks`n- Report on re-order levels`n- Detailed product sale/purchase history"

# Update quantity (hours) for the line item
$ws.Range("E14").Value2 = 68

# Update discount amount
$ws.Range("H19").Value2 = 200

# Update "Previous payment as per Invoice #" label to reference invoice 1
$ws.Range("G21").Value2 = "Previous payment as per Invoice # 1 : "

# Update the view so the window shows rows starting at 13 with B24 selected
$ws.Range("B24").Select()
$excel.ActiveWindow.ScrollRow = 13

$wb.Save()
